$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.406999999999999
$ws.Range("A3").Value = -21.995
$ws.Range("A14").Value = -21.659
$ws.Range("A21").Value = -19.884
$ws.Range("A23").Value = -20.317
$ws.Range("A25").Value = -21.759
$ws.Range("D25").Value = -8.040000000000001
$ws.Range("A26").Value = -21.351
$ws.Range("D27").Value = -8.653
$ws.Range("A29").Value = -21.297
$ws.Range("D31").Value = -8.161
$ws.Range("D39").Value = -7.502000000000001
$ws.Range("D48").Value = -7.475
$ws.Range("D51").Value = -8.259
$ws.Range("D52").Value = -8.199999999999999
$ws.Range("A53").Value = -22.01
$ws.Range("D55").Value = -8.114999999999998
$ws.Range("D56").Value = -8.501999999999999
$ws.Range("A57").Value = -22.473
$ws.Range("D57").Value = -8.084
$ws.Range("A59").Value = -22.5
$ws.Range("A69").Value = -21.602
$ws.Range("D73").Value = -8.004000000000001
$ws.Range("A79").Value = -20.973
$ws.Range("A83").Value = -22.038
$ws.Range("D89").Value = -6.823
$ws.Range("D90").Value = -7.601999999999999
$ws.Range("A91").Value = -21.509
$ws.Range("D92").Value = -6.621
$ws.Range("A93").Value = -21.558
